$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Unprotect the sheet so values can be updated (sheet was protected with a password)
$ws.Unprotect("lido")

# Update the confidential disclaimer text (shared string) with the new date
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-07-08 for illustrative purposes only and are subject to change."

# Update the weight / percent-change figures
$ws.Range("D2").Value = 0.8468146945055884
$ws.Range("E2").Value = -0.01287066246056778

$ws.Range("D3").Value = 0.1531853054944116
$ws.Range("E3").Value = -0.02083333333333326

$ws.Range("E4").Value = -0.01409042663076387

# Re-protect the sheet to restore original protection settings
$ws.Protect("lido", $false, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false, $true, $true)

$wb.Save()
